$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.456.72'
$ws.Range("E2").Value = '  +3.35%  '

$ws.Range("D3").Value = '2.634.98'
$ws.Range("E3").Value = '  +1.06%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.70%  '

$ws.Range("E7").Value = '  -0.39%  '

$ws.Range("E8").Value = '  +4.15%  '

$ws.Range("D9").Value = '2.660.88'
$ws.Range("E9").Value = '  +1.89%  '

$ws.Range("E10").Value = '  +0.43%  '

$ws.Range("E11").Value = '  +5.09%  '

$ws.Range("E12").Value = '  +6.82%  '

$ws.Range("E13").Value = '  +2.98%  '

$ws.Range("D14").Value = '3.111.40'
$ws.Range("E14").Value = '  +1.31%  '

$ws.Range("D15").Value = '60.476.24'
$ws.Range("E15").Value = '  +3.46%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.30%  '

$ws.Range("E17").Value = '  +4.89%  '

$ws.Range("D18").Value = '2.647.06'
$ws.Range("E18").Value = '  +1.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.99%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.440'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.59%  '

$ws.Range("E26").Value = '  +2.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.995'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.28%  '

$ws.Range("E28").Value = '  +4.80%  '

$ws.Range("D29").Value = '0.0₃0815'
$ws.Range("E29").Value = '  +11.08%  '

$ws.Range("E31").Value = '  +4.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '159.37'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.85%  '

$ws.Range("E34").Value = '  +2.05%  '

$ws.Range("E35").Value = '  +5.54%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.902'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.14%  '

$ws.Range("E37").Value = '  +4.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.895'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.79%  '

$ws.Range("E39").Value = '  +7.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.44'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '300.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.58%  '

$ws.Range("E42").Value = '  +1.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.994'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.65%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0986'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.69%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.604'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '129.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +15.40%  '

$ws.Range("E47").Value = '  +2.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.47%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("E50").Value = '  +4.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.69'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.20%  '

